# "Plan Iteration 3" Action Plan sheet update:
#  - Add due-date text ("Nth March") into column D for each activity row.
#  - Rename the #4 task from "UCD - Saving the game" to
#    "UCD - First move of the token".
#  - Consolidate the last two board-setup rows: "Putting up pieces" /
#    "Saving the game" become "Moving up pieces" (row 21) and
#    "Setting up barricade" (row 22, now owned by Mohammad/Tomal) - the
#    redundant trailing row is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Due-date column (D) for the use-case / diagram activities
$ws.Range("D9").Value  = "14th March"
$ws.Range("D11").Value = "15th March"
$ws.Range("D13").Value = "17th March"
$ws.Range("D15").Value = "17th March"
$ws.Range("D17").Value = "15th March"
$ws.Range("D19").Value = "17th March"

# Task #4 renamed
$ws.Range("B11").Value = "#4: UCD - First move of the token"

# Board-setup section (rows 20-22) gains due dates too, and the
# "Putting up pieces" / "Saving the game" rows are reworked
$ws.Range("D20").Value = "13th March"

$ws.Range("B21").Value = "Moving up pieces"
$ws.Range("D21").Value = "15th March"

$ws.Range("B22").Value = "Setting up barricade"
$ws.Range("C22").Value = "Mohammad/Tomal"
$ws.Range("D22").Value = "16th March"

# The old trailing "Setting up barricade" row (23) is now redundant -
# its content was folded into row 22 above, so delete the whole row.
$ws.Rows("23:23").Delete()

# Leave the selection where the author left it when saving.
$ws.Range("C24").Select()
